# Revision restore (SAVE): cell C10 on the "Rules" sheet changes value 18 -> 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
